$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "['year', 'hour']"
$ws.Range("E2").Value = "['year', 'weekend', 'hour']"
$ws.Range("F2").Value = "['year', 'state', 'hour']"
$ws.Range("G2").Value = "['year', 'state', 'weekend', 'hour']"

$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.64)"
$ws.Range("D4").Value = "(0.6)"
$ws.Range("E4").Value = "(0.03)"
$ws.Range("F4").Value = "(0.2)"
$ws.Range("G4").Value = "(0.65)"

$ws.Range("B6").Value = "(0.0)"
$ws.Range("C6").Value = "(0.49)"
$ws.Range("D6").Value = "(0.51)"
$ws.Range("E6").Value = "(0.11)"
$ws.Range("F6").Value = "(0.51)"
$ws.Range("G6").Value = "(0.52)"
